$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Home_Page")

$ws.Range("A4").Value = "name"
$ws.Range("B4").Value = "AGYSITR"

$ws.Activate()
$ws.Range("A5").Select()
